# Pad the hour-less "raw!E" time strings with a leading zero so that
# single-digit hours (e.g. "9:00 PM") render consistently as "09:00 PM".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("transformed")

for ($row = 2; $row -le 23; $row++) {
    $ws.Range("B$row").Formula = '="0"&raw!E' + $row
}

# Leave the cursor on C19, matching where editing wrapped up.
$ws.Range("C19").Select() | Out-Null
